# "Generate Report for Handoff"
# The two tracked source files moved from "Handed back: in sync with en-US"
# to "Ready for handoff" for a fresh localization package. Update the
# Overview roll-up sheet plus the per-locale (zh-cn / de-de) detail sheets:
# new source-file identifiers, new handoff artifact names/timestamps, a
# cleared (not-yet-handed-back) Target/Handback pair, and refreshed
# hyperlinks everywhere the old identifiers used to show up.

$wb = $excel.ActiveWorkbook

# ---- old/new identifiers -------------------------------------------------
$oldUuid1 = "7d23684b-c227-4902-a598-92575ecc296d"
$oldHash1 = "92c79aae8c947b1a33993a8ecbcd6485dc405b01"
$oldUuid2 = "a865a3d9-612c-4cf4-ba43-b184ebb865f4"
$oldHash2 = "aa0b67468c58d6b3a9d67508302dd92113c61cd5"

$newUuid1 = "aa846359-e5b7-4b1a-992e-45eab8e66c07"
$newHash1 = "6048720f6a2b7d4f85484e10c32f8fbb6ee781b6"
$newUuid2 = "ffff54f9520a-8ea1-4e33-ba88-650474803040"

$newStatus = "Ready for handoff"

$newFile1 = $newUuid1 + ".md"
$newFile2 = $newUuid2 + ".md"

$newXlfZh = $newUuid1 + "." + $newHash1 + ".zh-cn.xlf"
$newXlfDe = $newUuid1 + "." + $newHash1 + ".de-de.xlf"

# ===========================================================================
# Overview sheet
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = "2016-44-11 22:44:30"

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-44-11 22:44:30"

# Hyperlinks: recreate A2/A3 to point at the new file names. The engine's
# hyperlink-mutation bridge only appends new entries when touched in place,
# so drop the whole collection for this sheet and rebuild it in the
# original order (A2 = rId2, A3 = rId3).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile1,
    [Type]::Missing, [Type]::Missing, $newFile1)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile2,
    [Type]::Missing, [Type]::Missing, $newFile2)

# ===========================================================================
# zh-cn sheet
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = "2016-03-11 22:44:27"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = "2016-03-11 22:44:27"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"

# The handoff hasn't produced a target/handback pair yet, so those two
# columns go blank for both rows.
$wsZh.Range("F2:G3").Clear()

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile1,
    [Type]::Missing, [Type]::Missing, $newFile1)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile1,
    [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6cd70f055be43ffe029cfff93f11514c8447e4a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newXlfZh,
    [Type]::Missing, [Type]::Missing, $newXlfZh)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile2,
    [Type]::Missing, [Type]::Missing, $newFile2)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile2,
    [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6cd70f055be43ffe029cfff93f11514c8447e4a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newXlfZh,
    [Type]::Missing, [Type]::Missing, $newXlfZh)

# ===========================================================================
# de-de sheet
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = "2016-03-11 22:44:30"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = "2016-03-11 22:44:30"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"

$wsDe.Range("F2:G3").Clear()

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile1,
    [Type]::Missing, [Type]::Missing, $newFile1)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile1,
    [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f15faa07752f88997fb25e8dd5ebcd3a04afee2f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newXlfDe,
    [Type]::Missing, [Type]::Missing, $newXlfDe)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile2,
    [Type]::Missing, [Type]::Missing, $newFile2)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/fd125f18fdd242c53b2a9eba76cfe53859678f83/e2e/" + $newFile2,
    [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f15faa07752f88997fb25e8dd5ebcd3a04afee2f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newXlfDe,
    [Type]::Missing, [Type]::Missing, $newXlfDe)
